# This script updates the TestReport worksheet to reflect a re-run of the
# test suite a few minutes later, where the final test case ("Verify Home
# Page Loads Successfully") now PASSES instead of FAILS.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update execution times for rows 2-9 (shift ~ +5m37s later re-run).
$ws.Range("D2").Value = "01/04/2025 01:36:14 PM"
$ws.Range("D3").Value = "01/04/2025 01:36:16 PM"
$ws.Range("D4").Value = "01/04/2025 01:36:17 PM"
$ws.Range("D5").Value = "01/04/2025 01:36:17 PM"
$ws.Range("D6").Value = "01/04/2025 01:36:20 PM"
$ws.Range("D7").Value = "01/04/2025 01:36:32 PM"
$ws.Range("D8").Value = "01/04/2025 01:36:34 PM"
$ws.Range("D9").Value = "01/04/2025 01:36:35 PM"

# Row 10 ("Verify Home Page Loads Successfully") now passes.
$ws.Range("C10").Value = "PASSED"
$ws.Range("D10").Value = "01/04/2025 01:36:35 PM"
$ws.Range("E10").Value = "Test executed successfully."
